$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.264634507416061
$ws.Range("D2").Value = 4.131892350614923
$ws.Range("E2").Value = 12.94353012288232
$ws.Range("F2").Value = 19.84064728920296
$ws.Range("G2").Value = 21.14119518385853
$ws.Range("H2").Value = 12.34772997530987
$ws.Range("I2").Value = 19.00773875537485
$ws.Range("K2").Value = 11.21451152523676
$ws.Range("M2").Value = 13.14410739126466
$ws.Range("O2").Value = 17.77530969758811
$ws.Range("B3").Value = 6.121846432533737
$ws.Range("D3").Value = 4.04615649401988
$ws.Range("E3").Value = 12.83110169929033
$ws.Range("F3").Value = 19.85857181684185
$ws.Range("G3").Value = 21.15806676827568
$ws.Range("H3").Value = 12.39362983825623
$ws.Range("I3").Value = 19.12678738975166
$ws.Range("K3").Value = 10.83460530348953
$ws.Range("M3").Value = 12.86686642439832
$ws.Range("O3").Value = 17.8424945560641
$ws.Range("B4").Value = 6.033133252159688
$ws.Range("D4").Value = 3.991964515154414
$ws.Range("E4").Value = 12.76689868771772
$ws.Range("F4").Value = 19.87606788453306
$ws.Range("G4").Value = 21.17789419210035
$ws.Range("H4").Value = 12.42407114179932
$ws.Range("I4").Value = 19.20361581150247
$ws.Range("K4").Value = 10.59285575717776
$ws.Range("M4").Value = 12.69605493571375
$ws.Range("O4").Value = 17.88843762545484
$ws.Range("B5").Value = 5.996769713533097
$ws.Range("D5").Value = 3.969508424055047
$ws.Range("E5").Value = 12.741977943029
$ws.Range("F5").Value = 19.8848261948227
$ws.Range("G5").Value = 21.18834493267646
$ws.Range("H5").Value = 12.43704378126362
$ws.Range("I5").Value = 19.23586517630304
$ws.Range("K5").Value = 10.49230843079914
$ws.Range("M5").Value = 12.62639935856384
$ws.Range("O5").Value = 17.90833558631832
$ws.Range("B6").Value = 5.990720246376793
$ws.Range("D6").Value = 3.965757637923495
$ws.Range("E6").Value = 12.73791568498308
$ws.Range("F6").Value = 19.88637873988914
$ws.Range("G6").Value = 21.19022313018539
$ws.Range("H6").Value = 12.43923214632409
$ws.Range("I6").Value = 19.24127708189958
$ws.Range("K6").Value = 10.4754930419136
$ws.Range("M6").Value = 12.61483309972584
$ws.Range("O6").Value = 17.91171052635462
$ws.Range("B7").Value = 6.032643633009338
$ws.Range("D7").Value = 3.99166314883628
$ws.Range("E7").Value = 12.76655753180666
$ws.Range("F7").Value = 19.87617941382791
$ws.Range("G7").Value = 21.17802554965502
$ws.Range("H7").Value = 12.42424379782936
$ws.Range("I7").Value = 19.20404692340813
$ws.Range("K7").Value = 10.5915078282382
$ws.Range("M7").Value = 12.69511560070594
$ws.Range("O7").Value = 17.8887012208292
$ws.Range("B8").Value = 6.215646219874276
$ws.Range("D8").Value = 4.102662368725765
$ws.Range("E8").Value = 12.903780569888
$ws.Range("F8").Value = 19.84547880165471
$ws.Range("G8").Value = 21.1450421166094
$ws.Range("H8").Value = 12.36308724644606
$ws.Range("I8").Value = 19.04801375943824
$ws.Range("K8").Value = 11.08533857759361
$ws.Range("M8").Value = 13.04869421999375
$ws.Range("O8").Value = 17.79749906668952
$ws.Range("B9").Value = 6.564194355655459
$ws.Range("D9").Value = 4.307324827677403
$ws.Range("E9").Value = 13.20979723977434
$ws.Range("F9").Value = 19.8368814642292
$ws.Range("G9").Value = 21.15581694296206
$ws.Range("H9").Value = 12.26109741763884
$ws.Range("I9").Value = 18.77152162338909
$ws.Range("K9").Value = 11.98248582101133
$ws.Range("M9").Value = 13.7332500756506
$ws.Range("O9").Value = 17.65604320459275
$ws.Range("B10").Value = 6.811396334619589
$ws.Range("D10").Value = 4.448882547854841
$ws.Range("E10").Value = 13.45509401860573
$ws.Range("F10").Value = 19.86210133134809
$ws.Range("G10").Value = 21.21001074021543
$ws.Range("H10").Value = 12.19712401969727
$ws.Range("I10").Value = 18.5861850829817
$ws.Range("K10").Value = 12.59361888708364
$ws.Range("M10").Value = 14.22549496576903
$ws.Range("O10").Value = 17.57513170484575
$ws.Range("B11").Value = 6.921434384991549
$ws.Range("D11").Value = 4.511202495093068
$ws.Range("E11").Value = 13.57065156759668
$ws.Range("F11").Value = 19.88041347897431
$ws.Range("G11").Value = 21.24471749901811
$ws.Range("H11").Value = 12.17040543073298
$ws.Range("I11").Value = 18.50569914206797
$ws.Range("K11").Value = 12.86043184396302
$ws.Range("M11").Value = 14.44606695064005
$ws.Range("O11").Value = 17.54336490194993
$ws.Range("B12").Value = 6.962718314322595
$ws.Range("D12").Value = 4.534490440146683
$ws.Range("E12").Value = 13.61493890849104
$ws.Range("F12").Value = 19.88832850148185
$ws.Range("G12").Value = 21.25930172779445
$ws.Range("H12").Value = 12.16063096148804
$ws.Range("I12").Value = 18.47576841374336
$ws.Range("K12").Value = 12.95980233086882
$ws.Range("M12").Value = 14.52902970618524
$ws.Range("O12").Value = 17.53206378356143
$ws.Range("B13").Value = 6.953844779315364
$ws.Range("D13").Value = 4.52948900201006
$ws.Range("E13").Value = 13.60537803737054
$ws.Range("F13").Value = 19.8865803065886
$ws.Range("G13").Value = 21.25609673652015
$ws.Range("H13").Value = 12.16272079169537
$ws.Range("I13").Value = 18.48219021741084
$ws.Range("K13").Value = 12.938476020863
$ws.Range("M13").Value = 14.5111884679479
$ws.Range("O13").Value = 17.53446523755523
$ws.Range("B14").Value = 6.924838762621061
$ws.Range("D14").Value = 4.513124716589158
$ws.Range("E14").Value = 13.57428477655003
$ws.Range("F14").Value = 19.88104502089215
$ws.Range("G14").Value = 21.24588849558279
$ws.Range("H14").Value = 12.16959439612515
$ws.Range("I14").Value = 18.50322576307103
$ws.Range("K14").Value = 12.86864079398958
$ws.Range("M14").Value = 14.45290407617168
$ws.Range("O14").Value = 17.54242053698226
$ws.Range("B15").Value = 6.907020517757842
$ws.Range("D15").Value = 4.503060219182498
$ws.Range("E15").Value = 13.55530675646009
$ws.Range("F15").Value = 19.87778209511895
$ws.Range("G15").Value = 21.23982321219322
$ws.Range("H15").Value = 12.17384939960339
$ws.Range("I15").Value = 18.51618187807187
$ws.Range("K15").Value = 12.82564622648353
$ws.Range("M15").Value = 14.4171275623737
$ws.Range("O15").Value = 17.54738833164175
$ws.Range("B16").Value = 6.804153191886318
$ws.Range("D16").Value = 4.444766911639637
$ws.Range("E16").Value = 13.44761833098899
$ws.Range("F16").Value = 19.86104199854554
$ws.Range("G16").Value = 21.2079445841588
$ws.Range("H16").Value = 12.19891814818537
$ws.Range("I16").Value = 18.59152178146254
$ws.Range("K16").Value = 12.57595181399754
$ws.Range("M16").Value = 14.21100556312232
$ws.Range("O16").Value = 17.57730946639424
$ws.Range("B17").Value = 6.740401001776503
$ws.Range("D17").Value = 4.408464860327353
$ws.Range("E17").Value = 13.38254190671179
$ws.Range("F17").Value = 19.85252258956897
$ws.Range("G17").Value = 21.19096098034027
$ws.Range("H17").Value = 12.21490782540812
$ws.Range("I17").Value = 18.63871820943881
$ws.Range("K17").Value = 12.41986298527292
$ws.Range("M17").Value = 14.08363842979778
$ws.Range("O17").Value = 17.59695859330823
$ws.Range("B18").Value = 6.703507157253566
$ws.Range("D18").Value = 4.387390240564622
$ws.Range("E18").Value = 13.34548806487018
$ws.Range("F18").Value = 19.84826647142188
$ws.Range("G18").Value = 21.18213919091676
$ws.Range("H18").Value = 12.22432895142918
$ws.Range("I18").Value = 18.66622447719833
$ws.Range("K18").Value = 12.32903399716869
$ws.Range("M18").Value = 14.01006916646083
$ws.Range("O18").Value = 17.60873446590209
$ws.Range("B19").Value = 6.690978042762355
$ws.Range("D19").Value = 4.380221713972219
$ws.Range("E19").Value = 13.33300820258758
$ws.Range("F19").Value = 19.84693610180167
$ws.Range("G19").Value = 21.17931496867947
$ws.Range("H19").Value = 12.22755728673805
$ws.Range("I19").Value = 18.67559955377197
$ws.Range("K19").Value = 12.29810221007331
$ws.Range("M19").Value = 13.98510900619853
$ws.Range("O19").Value = 17.61280291097978
$ws.Range("B20").Value = 6.747211128634531
$ws.Range("D20").Value = 4.412349516234089
$ws.Range("E20").Value = 13.38943075411808
$ws.Range("F20").Value = 19.85336285847737
$ws.Range("G20").Value = 21.19267095244193
$ws.Range("H20").Value = 12.21318247997376
$ws.Range("I20").Value = 18.63365681453567
$ws.Range("K20").Value = 12.43658809259265
$ws.Range("M20").Value = 14.09722967489209
$ws.Range("O20").Value = 17.59481780479058
$ws.Range("B21").Value = 6.933369270848866
$ws.Range("D21").Value = 4.517939848953294
$ws.Range("E21").Value = 13.58340363178538
$ws.Range("F21").Value = 19.88264428442918
$ws.Range("G21").Value = 21.2488478276024
$ws.Range("H21").Value = 12.167566130776
$ws.Range("I21").Value = 18.49703227218867
$ws.Range("K21").Value = 12.889198712859
$ws.Range("M21").Value = 14.47003951543437
$ws.Range("O21").Value = 17.54006408048374
$ws.Range("B22").Value = 7.052774269417377
$ws.Range("D22").Value = 4.585129672783888
$ws.Range("E22").Value = 13.71323440550257
$ws.Range("F22").Value = 19.90749489183849
$ws.Range("G22").Value = 21.29396113711012
$ws.Range("H22").Value = 12.13975437672137
$ws.Range("I22").Value = 18.41093067563459
$ws.Range("K22").Value = 13.17527537136428
$ws.Range("M22").Value = 14.71037649005302
$ws.Range("O22").Value = 17.50852581175639
$ws.Range("B23").Value = 6.989264109846149
$ws.Range("D23").Value = 4.549439675413259
$ws.Range("E23").Value = 13.64367583108475
$ws.Range("F23").Value = 19.89371015533005
$ws.Range("G23").Value = 21.26911694613496
$ws.Range("H23").Value = 12.15441473685059
$ws.Range("I23").Value = 18.45659358551428
$ws.Range("K23").Value = 13.02349790190765
$ws.Range("M23").Value = 14.58243290015783
$ws.Range("O23").Value = 17.52496872459795
$ws.Range("B24").Value = 6.744133020736041
$ws.Range("D24").Value = 4.410593897999987
$ws.Range("E24").Value = 13.38631518328172
$ws.Range("F24").Value = 19.8529809735237
$ws.Range("G24").Value = 21.19189493853596
$ws.Range("H24").Value = 12.21396179707738
$ws.Range("I24").Value = 18.63594391094526
$ws.Range("K24").Value = 12.42903006688163
$ws.Range("M24").Value = 14.0910861388792
$ws.Range("O24").Value = 17.5957841624577
$ws.Range("B25").Value = 6.471284878630271
$ws.Range("D25").Value = 4.253445281797457
$ws.Range("E25").Value = 13.12326081092765
$ws.Range("F25").Value = 19.83366989903309
$ws.Range("G25").Value = 21.14478502770749
$ws.Range("H25").Value = 12.28676560981872
$ws.Range("I25").Value = 18.8431812844154
$ws.Range("K25").Value = 11.74791898796407
$ws.Range("M25").Value = 13.54956928522681
$ws.Range("O25").Value = 17.69028428399887
